$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was renamed from "Property1" to "DataNode" as part of unifying
# the DataNode / DataTable / Entity naming convention across the config
# workbooks.
$ws.Name = "DataNode"
